$d = $word.ActiveDocument
$r = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:tbl><w:tblPr><w:tblpPr w:leftFromText="141" w:rightFromText="141" w:vertAnchor="text" w:horzAnchor="margin" w:tblpY="654"/><w:tblW w:w="10456" w:type="dxa"/><w:tblBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:insideH w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:insideV w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tblBorders><w:tblLook w:val="01E0"/></w:tblPr><w:tblGrid><w:gridCol w:w="6629"/><w:gridCol w:w="3827"/></w:tblGrid><w:tr><w:trPr><w:trHeight w:val="351"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="10456" w:type="dxa"/><w:gridSpan w:val="2"/><w:shd w:val="clear" w:color="auto" w:fill="4BACC6"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Encabezado"/><w:tabs><w:tab w:val="left" w:pos="2822"/></w:tabs><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>Tipo de Revisión: de Análisis de R</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>iesgo</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="347"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="6629" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Encabezado"/><w:tabs><w:tab w:val="left" w:pos="2108"/><w:tab w:val="right" w:pos="9602"/></w:tabs><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>Fecha de la revisión:</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3827" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Encabezado"/><w:tabs><w:tab w:val="left" w:pos="2108"/><w:tab w:val="right" w:pos="9602"/></w:tabs><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>Hora:</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="358"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="10456" w:type="dxa"/><w:gridSpan w:val="2"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Encabezado"/><w:tabs><w:tab w:val="left" w:pos="2108"/><w:tab w:val="left" w:pos="3802"/></w:tabs><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>Nombre del Riesgo:</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="358"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="10456" w:type="dxa"/><w:gridSpan w:val="2"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Encabezado"/><w:tabs><w:tab w:val="left" w:pos="2108"/><w:tab w:val="left" w:pos="3802"/></w:tabs><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>Encargado:</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="10456" w:type="dxa"/><w:gridSpan w:val="2"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Encabezado"/><w:tabs><w:tab w:val="left" w:pos="2108"/><w:tab w:val="right" w:pos="9602"/></w:tabs><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>Preguntas de comprobación:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Encabezado"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="2108"/><w:tab w:val="right" w:pos="9602"/></w:tabs><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:left="714" w:hanging="357"/><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>¿Se encuentra el elemento en revisión actualizado?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Encabezado"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="2108"/><w:tab w:val="right" w:pos="9602"/></w:tabs><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:left="714" w:hanging="357"/><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>¿Se realiza el plan de acción del riesgo?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Encabezado"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="2108"/><w:tab w:val="right" w:pos="9602"/></w:tabs><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:left="714" w:hanging="357"/><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>¿El Riesgo considerado</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t xml:space="preserve"> se puede llegar a concretar?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Encabezado"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="2108"/><w:tab w:val="right" w:pos="9602"/></w:tabs><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:left="714" w:hanging="357"/><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>¿Se llego a concretar</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t xml:space="preserve"> algun</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>o de los disparadores del Riesgo?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Encabezado"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="2108"/><w:tab w:val="right" w:pos="9602"/></w:tabs><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:left="714" w:hanging="357"/><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t xml:space="preserve">¿Es </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>necesario realizar</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t xml:space="preserve"> el plan de co</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>ntingencias del riesgo considerado</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Encabezado"/><w:tabs><w:tab w:val="left" w:pos="2108"/><w:tab w:val="right" w:pos="9602"/></w:tabs><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:left="714"/><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Encabezado"/><w:tabs><w:tab w:val="left" w:pos="2108"/><w:tab w:val="right" w:pos="9602"/></w:tabs><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>Resultado Obtenido:__________________________________________________</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Encabezado"/><w:tabs><w:tab w:val="left" w:pos="2108"/><w:tab w:val="right" w:pos="9602"/></w:tabs><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Encabezado"/><w:tabs><w:tab w:val="left" w:pos="2108"/><w:tab w:val="right" w:pos="9602"/></w:tabs><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>Otras Observaciones:_________________________________________________</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Encabezado"/><w:tabs><w:tab w:val="left" w:pos="2108"/><w:tab w:val="right" w:pos="9602"/></w:tabs><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr></w:pPr></w:p></w:tc></w:tr></w:tbl><w:p/><w:p/><w:p/><w:p/><w:p/><w:p/><w:p/><w:p/><w:p/><w:p/><w:p/><w:tbl><w:tblPr><w:tblpPr w:leftFromText="141" w:rightFromText="141" w:vertAnchor="text" w:horzAnchor="page" w:tblpX="673" w:tblpY="-755"/><w:tblW w:w="10598" w:type="dxa"/><w:tblBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:insideH w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:insideV w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tblBorders><w:tblLook w:val="01E0"/></w:tblPr><w:tblGrid><w:gridCol w:w="5990"/><w:gridCol w:w="4608"/></w:tblGrid><w:tr><w:trPr><w:trHeight w:val="351"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="10598" w:type="dxa"/><w:gridSpan w:val="2"/><w:shd w:val="clear" w:color="auto" w:fill="4BACC6"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Encabezado"/><w:tabs><w:tab w:val="left" w:pos="2822"/></w:tabs><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:lastRenderedPageBreak/><w:t>Tipo de Revisión: de Análisis de R</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>iesgo</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="347"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="5990" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Encabezado"/><w:tabs><w:tab w:val="left" w:pos="2108"/><w:tab w:val="right" w:pos="9602"/></w:tabs><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>Fecha de la revisión:</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t xml:space="preserve"> 17/05/2013</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4608" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Encabezado"/><w:tabs><w:tab w:val="left" w:pos="2108"/><w:tab w:val="right" w:pos="9602"/></w:tabs><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>Hora:</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>17:00</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="358"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="10598" w:type="dxa"/><w:gridSpan w:val="2"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Encabezado"/><w:tabs><w:tab w:val="left" w:pos="2108"/><w:tab w:val="left" w:pos="3802"/></w:tabs><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>Nombre del Riesgo</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>:</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> RI- 13 Planificación incorrecta del cronograma de desarrollo del proyecto.</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="358"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="10598" w:type="dxa"/><w:gridSpan w:val="2"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Encabezado"/><w:tabs><w:tab w:val="left" w:pos="2108"/><w:tab w:val="left" w:pos="3802"/></w:tabs><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>Encargado:</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t xml:space="preserve"></w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>Lizza</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t xml:space="preserve"> Lorena López </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>Maciel</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>.</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="10598" w:type="dxa"/><w:gridSpan w:val="2"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Encabezado"/><w:tabs><w:tab w:val="left" w:pos="2108"/><w:tab w:val="right" w:pos="9602"/></w:tabs><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>Preguntas de comprobación:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Encabezado"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="2108"/><w:tab w:val="right" w:pos="9602"/></w:tabs><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:left="714" w:hanging="357"/><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>¿Se encuentra el elemento en revisión actualizado?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Encabezado"/><w:tabs><w:tab w:val="left" w:pos="2108"/><w:tab w:val="right" w:pos="9602"/></w:tabs><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:left="714"/><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>Si.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Encabezado"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="2108"/><w:tab w:val="right" w:pos="9602"/></w:tabs><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:left="714" w:hanging="357"/><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>¿Se realiza el plan de acción del riesgo?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Encabezado"/><w:tabs><w:tab w:val="left" w:pos="2108"/><w:tab w:val="right" w:pos="9602"/></w:tabs><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:left="714"/><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>Se llevo a cabo el plan de acción, ya que se fijo los objetivos a cumplir en cierto tiempo y los integrantes van informando o actualizando lo que van haciendo.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Encabezado"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="2108"/><w:tab w:val="right" w:pos="9602"/></w:tabs><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:left="714" w:hanging="357"/><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>¿El Riesgo considerado</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t xml:space="preserve"> se puede llegar a concretar?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Encabezado"/><w:tabs><w:tab w:val="left" w:pos="2108"/><w:tab w:val="right" w:pos="9602"/></w:tabs><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:left="714"/><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>Si.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Encabezado"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="2108"/><w:tab w:val="right" w:pos="9602"/></w:tabs><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:left="714" w:hanging="357"/><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>¿Se llego a concretar</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t xml:space="preserve"> algun</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>o de los disparadores del Riesgo?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Encabezado"/><w:tabs><w:tab w:val="left" w:pos="2108"/><w:tab w:val="right" w:pos="9602"/></w:tabs><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:left="714"/><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>Si, el tiempo de desarrollo no es suficiente.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Encabezado"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="2108"/><w:tab w:val="right" w:pos="9602"/></w:tabs><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:left="714" w:hanging="357"/><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t xml:space="preserve">¿Es </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>necesario realizar</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t xml:space="preserve"> el plan de co</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>ntingencias del riesgo considerado</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Encabezado"/><w:tabs><w:tab w:val="left" w:pos="2108"/><w:tab w:val="right" w:pos="9602"/></w:tabs><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:left="714"/><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t xml:space="preserve">Si, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>se va realizando ajustes de a qué actividad hay que darle prioridad.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Encabezado"/><w:tabs><w:tab w:val="left" w:pos="2108"/><w:tab w:val="right" w:pos="9602"/></w:tabs><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:left="714"/><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Encabezado"/><w:tabs><w:tab w:val="left" w:pos="2108"/><w:tab w:val="right" w:pos="9602"/></w:tabs><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>Resul</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>tado Obtenido: Este riesgo en cierta medida se va cumpliendo, ya que ciertas actividades aún no han sido concluidas como se estableció.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t xml:space="preserve"> Pero, se va ajustando las actividades y viendo a que tarea se debe dar prioridad.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Encabezado"/><w:tabs><w:tab w:val="left" w:pos="2108"/><w:tab w:val="right" w:pos="9602"/></w:tabs><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Encabezado"/><w:tabs><w:tab w:val="left" w:pos="2108"/><w:tab w:val="right" w:pos="9602"/></w:tabs><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr><w:t>Otras Observaciones:_________________________________________________</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Encabezado"/><w:tabs><w:tab w:val="left" w:pos="2108"/><w:tab w:val="right" w:pos="9602"/></w:tabs><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="SimHei" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES_tradnl" w:eastAsia="ar-SA"/></w:rPr></w:pPr></w:p></w:tc></w:tr></w:tbl><w:p/><w:p/><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="2700"/></w:tabs></w:pPr><w:r><w:tab/></w:r></w:p><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="2700"/></w:tabs></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)
Write-Host "Tables after: " $d.Tables.Count
Write-Host "Paragraphs after: " $d.Paragraphs.Count
